$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2052980132450331
$ws.Range("C2").Value = 0.543046357615894
$ws.Range("J2").Value = 0.01655629139072848
$ws.Range("P2").Value = 0.1556291390728477
$ws.Range("S2").Value = 0.07947019867549669
$ws.Range("C3").Value = 0.02906976744186046
$ws.Range("J3").Value = 0.02325581395348837
$ws.Range("P3").Value = 0.7558139534883721
$ws.Range("S3").Value = 0.1918604651162791
$ws.Range("J4").Value = 0.03125
$ws.Range("P4").Value = 0.625
$ws.Range("S4").Value = 0.34375
$ws.Range("B6").Value = 0.0776255707762557
$ws.Range("D6").Value = 0.0045662100456621
$ws.Range("F6").Value = 0.0730593607305936
$ws.Range("J6").Value = 0.2100456621004566
$ws.Range("Q6").Value = 0.1780821917808219
$ws.Range("R6").Value = 0.0776255707762557
$ws.Range("S6").Value = 0.3789954337899543
$ws.Range("B7").Value = 0.09223300970873786
$ws.Range("D7").Value = 0.009708737864077669
$ws.Range("E7").Value = 0.004854368932038835
$ws.Range("F7").Value = 0.05339805825242718
$ws.Range("J7").Value = 0.1601941747572816
$ws.Range("O7").Value = 0.01456310679611651
$ws.Range("Q7").Value = 0.1844660194174757
$ws.Range("R7").Value = 0.07766990291262135
$ws.Range("S7").Value = 0.4029126213592233
$ws.Range("B8").Value = 0.1113636363636364
$ws.Range("D8").Value = 0.01818181818181818
$ws.Range("E8").Value = 0.002272727272727273
$ws.Range("F8").Value = 0.04318181818181818
$ws.Range("J8").Value = 0.1272727272727273
$ws.Range("O8").Value = 0.01818181818181818
$ws.Range("Q8").Value = 0.1454545454545454
$ws.Range("R8").Value = 0.1159090909090909
$ws.Range("S8").Value = 0.4181818181818182
$ws.Range("B9").Value = 0.109452736318408
$ws.Range("D9").Value = 0.01492537313432836
$ws.Range("F9").Value = 0.06965174129353234
$ws.Range("J9").Value = 0.154228855721393
$ws.Range("O9").Value = 0.009950248756218905
$ws.Range("Q9").Value = 0.1741293532338309
$ws.Range("R9").Value = 0.08955223880597014
$ws.Range("S9").Value = 0.3781094527363184
$ws.Range("B10").Value = 0.107229894394801
$ws.Range("D10").Value = 0.01543460601137287
$ws.Range("F10").Value = 0.06986190089358245
$ws.Range("J10").Value = 0.1121039805036556
$ws.Range("O10").Value = 0.01299756295694557
$ws.Range("Q10").Value = 0.2030869212022746
$ws.Range("R10").Value = 0.08692120227457352
$ws.Range("S10").Value = 0.3923639317627945
$ws.Range("G11").Value = 0.1508196721311476
$ws.Range("J11").Value = 0.07868852459016394
$ws.Range("K11").Value = 0.2131147540983606
$ws.Range("L11").Value = 0.5409836065573771
$ws.Range("S11").Value = 0.01639344262295082
$ws.Range("G12").Value = 0.7705882352941177
$ws.Range("J12").Value = 0.2
$ws.Range("L12").Value = 0.01176470588235294
$ws.Range("S12").Value = 0.01764705882352941
$ws.Range("G13").Value = 0.7021276595744681
$ws.Range("J13").Value = 0.2553191489361702
$ws.Range("S13").Value = 0.0425531914893617
$ws.Range("F15").Value = 0.03061224489795918
$ws.Range("H15").Value = 0.2142857142857143
$ws.Range("I15").Value = 0.0663265306122449
$ws.Range("J15").Value = 0.3316326530612245
$ws.Range("K15").Value = 0.07142857142857142
$ws.Range("M15").Value = 0.02040816326530612
$ws.Range("O15").Value = 0.04591836734693878
$ws.Range("S15").Value = 0.2193877551020408
$ws.Range("F16").Value = 0.01595744680851064
$ws.Range("H16").Value = 0.1542553191489362
$ws.Range("I16").Value = 0.101063829787234
$ws.Range("J16").Value = 0.3776595744680851
$ws.Range("K16").Value = 0.1276595744680851
$ws.Range("M16").Value = 0.02659574468085106
$ws.Range("N16").Value = 0.005319148936170213
$ws.Range("O16").Value = 0.03723404255319149
$ws.Range("S16").Value = 0.1542553191489362
$ws.Range("F17").Value = 0.01418439716312057
$ws.Range("H17").Value = 0.182033096926714
$ws.Range("I17").Value = 0.1111111111111111
$ws.Range("J17").Value = 0.3853427895981087
$ws.Range("K17").Value = 0.08983451536643026
$ws.Range("M17").Value = 0.01891252955082742
$ws.Range("O17").Value = 0.06619385342789598
$ws.Range("S17").Value = 0.132387706855792
$ws.Range("F18").Value = 0.01456310679611651
$ws.Range("H18").Value = 0.1504854368932039
$ws.Range("I18").Value = 0.09223300970873786
$ws.Range("J18").Value = 0.4271844660194175
$ws.Range("K18").Value = 0.0825242718446602
$ws.Range("M18").Value = 0.009708737864077669
$ws.Range("N18").Value = 0.004854368932038835
$ws.Range("O18").Value = 0.07766990291262135
$ws.Range("S18").Value = 0.1407766990291262
$ws.Range("F19").Value = 0.01857585139318885
$ws.Range("H19").Value = 0.2043343653250774
$ws.Range("I19").Value = 0.08126934984520123
$ws.Range("J19").Value = 0.3653250773993808
$ws.Range("K19").Value = 0.1130030959752322
$ws.Range("M19").Value = 0.02476780185758514
$ws.Range("O19").Value = 0.0673374613003096
$ws.Range("S19").Value = 0.1253869969040248
